# Update the script for 0.7.0
$wb = $excel.ActiveWorkbook

# --- Sheet4 "Tribuanl Quests" : Karrod notes ---
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet5 "Bloodmoon" : Stones table + Karrod/BM startup notes ---
$ws5 = $wb.Worksheets.Item(5)

# Header row
$ws5.Range("B2").Value = "Stones:"
$ws5.Range("C2").Value = "Quest ID"
$ws5.Range("D2").Value = "Completion index"

# Beast row
$ws5.Range("B3").Value = "Beast"
$ws5.Range("C3").Value = "BM_Beasts"

$ws5.Range("F2").Value = "BM_Stones index"
$ws5.Range("E2").Value = "Done index"

$ws5.Range("B4").Value = "Earth"

$ws5.Range("G2").Value = "Particles"

$ws5.Range("G3").Value = "Act_BM_Beast_parts"

$ws5.Range("C4").Value = "BM_Earth"

$ws5.Range("G4").Value = "Act_BM_Earth_parts"

$ws5.Range("B5").Value = "Sun"
$ws5.Range("B6").Value = "Trees"
$ws5.Range("B7").Value = "Water"
$ws5.Range("B8").Value = "Wind"

$ws5.Range("G5").Value = "Act_BM_Sun_parts"

$ws5.Range("C5").Value = "BM_Sun"

$ws5.Range("G6").Value = "Act_BM_Tree_parts"

$ws5.Range("C6").Value = "BM_Trees"

$ws5.Range("G7").Value = "Act_BM_Water_parts"

$ws5.Range("C7").Value = "BM_Water"

$ws5.Range("C8").Value = "BM_Wind"

$ws5.Range("G8").Value = "Act_BM_Wind_parts"

# Numeric columns for each stone row
$ws5.Range("D3").Value = 50
$ws5.Range("E3").Value = 100
$ws5.Range("F3").Value = 66

$ws5.Range("D4").Value = 50
$ws5.Range("E4").Value = 100
$ws5.Range("F4").Value = 62

$ws5.Range("D5").Value = 40
$ws5.Range("E5").Value = 100
$ws5.Range("F5").Value = 68

$ws5.Range("D6").Value = 40
$ws5.Range("E6").Value = 100
$ws5.Range("F6").Value = 64

$ws5.Range("D7").Value = 70
$ws5.Range("E7").Value = 100
$ws5.Range("F7").Value = 60

$ws5.Range("D8").Value = 50
$ws5.Range("E8").Value = 100
$ws5.Range("F8").Value = 70

# Column widths (match the workbook's standard auto-fit look used on other sheets)
$ws5.Range("C:C").ColumnWidth = 17.7
$ws5.Range("D:D").ColumnWidth = 17.7
$ws5.Range("E:E").ColumnWidth = 17.5
$ws5.Range("F:F").ColumnWidth = 17.5
$ws5.Range("G:G").ColumnWidth = 17.5

# Notes beneath the table
$ws5.Range("B10").Value = "for each stone, set stones to ( stones ) +1"
$ws5.Range("B11").Value = "additionally, set doOnce to 1 on stones if quest index > 10, 2 if quest index = 100"
$ws5.Range("B12").Value = "Act_BM_stone_<thing>_01"

# --- Sheet4 "Tribuanl Quests" : Karrod notes ---
$ws4.Range("A3").Value = "Karrod:"
$ws4.Range("A4").Value = "if TR_Champion == 20, set KarrodFightStart to 1, start karrodMovement script"
$ws4.Range("A5").Value = "if TR_Champion == 50, set KarrodBribe to 1"
$ws4.Range("A7").Value = "TR_Blade == 60 then start bladefixScript"
$ws4.Range("A7").Select()

# --- back to Sheet5 for the final note ---
$ws5.Range("B14").Value = "Start BMStartUpScript on empty server"

# --- Sheet1 "Startup" : move selection, this sheet is no longer the active tab ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("M41").Select()

# --- Activate Bloodmoon sheet (sheet index 5, 0-based activeTab = 4) and set its selection ---
$ws5.Activate()
$ws5.Range("F17").Select()
